$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Drop existing hyperlinks up front; this engine does not re-anchor
#     hyperlink ranges when columns are inserted, so we rebuild them all
#     (at their correct, post-insert addresses) once the layout is final.
$ws.Hyperlinks.Delete()

# --- Insert a new column before D. This pushes the old "phone" column
#     (D) to E, matching the target layout (A Name, B Username, C Email,
#     D <new>, E Phone, F <new>).
$ws.Range("D1").EntireColumn.Insert()

# --- Row 6 (Pierre-Luc Buhler) gets a second email + a new home phone.
#     The old university email moves from C6 into the newly freed D6,
#     and the new personal email takes its place in C6.
$oldEmail = $ws.Range("C6").Text

# Give D6 the same look as the other wrapped hyperlink cells (copy format
# from E6, which still carries the original C/D-column style).
$ws.Range("E6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D6").Value = $oldEmail

$ws.Range("C6").Value = "plbuhler01@gmail.com"
$ws.Range("C6").Style = "Hyperlink"

# New home-phone cell, plain formatting (no hyperlink).
$ws.Range("F6").Value = "(maison) 418-847-7722"

# --- Column widths per final layout.
$ws.Columns("C").ColumnWidth = 29.5703125
$ws.Columns("D").ColumnWidth = 28.140625
$ws.Columns("E").ColumnWidth = 15.28515625
$ws.Columns("F").ColumnWidth = 9.140625

# --- Rebuild all hyperlinks at their final addresses.
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:francis.valois1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E1"), "tel:581-307-0691", "", "tel:581-307-0691")

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:danthib76@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E2"), "tel:%28581%29 997-6656", "", "tel:%28581%29 997-6656")

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:g.oliviersylvain@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "tel:418-456-3734", "", "tel:418-456-3734")

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:diane.fournier.3@ulaval.ca")
$ws.Hyperlinks.Add($ws.Range("E4"), "tel:418-455-1155", "", "tel:418-455-1155")

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:philippe.bourdages@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "tel:418-563-9854", "", "tel:418-563-9854")

$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:pierre-luc.buhler.1@ulaval.ca")
$ws.Hyperlinks.Add($ws.Range("E6"), "tel:581-777-0237", "", "tel:581-777-0237")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:plbuhler01@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:emile.arsenault.1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E7"), "tel:418-931-8720", "", "tel:418-931-8720")

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:imanemouhtij@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("E8"), "tel:418-271-7831", "", "tel:418-271-7831")

# --- Selection cosmetics from the diff.
$ws.Range("D15").Select()

Write-Output "done"
